$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 14917424
$ws.Range("I19").Value = 11739891
$ws.Range("J19").Value = 20001476
$ws.Range("K19").Value = 11739891
$ws.Range("L19").Value = 20001476
$ws.Range("M19").Value = -11739716
$ws.Range("N19").Value = -20001826
$ws.Range("H86").Value = 1830.8462
$ws.Range("I86").Value = 1818.4546
$ws.Range("J86").Value = 1899
$ws.Range("K86").Value = 1818.4546
$ws.Range("L86").Value = 1899
$ws.Range("M86").Value = -695.4546
$ws.Range("N86").Value = -4145
$ws.Range("H88").Value = 5523
$ws.Range("I88").Value = 1040.6
$ws.Range("J88").Value = 11126
$ws.Range("K88").Value = 1040.6
$ws.Range("L88").Value = 11126
$ws.Range("M88").Value = -634.5999999999999
$ws.Range("N88").Value = -11938
$ws.Range("H89").Value = 1830.8462
$ws.Range("I89").Value = 1818.4546
$ws.Range("J89").Value = 1899
$ws.Range("K89").Value = 9092.273000000001
$ws.Range("L89").Value = 9495
$ws.Range("M89").Value = -3476.273000000001
$ws.Range("N89").Value = -20727
$ws.Range("H91").Value = 5523
$ws.Range("I91").Value = 1040.6
$ws.Range("J91").Value = 11126
$ws.Range("K91").Value = 1040.6
$ws.Range("L91").Value = 11126
$ws.Range("M91").Value = 363.4000000000001
$ws.Range("N91").Value = -13934
$ws.Range("H94").Value = 4240.7144
$ws.Range("I94").Value = 4240.7144
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4240.7144
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3789.7144
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 1238.75
$ws.Range("I100").Value = 1182
$ws.Range("J100").Value = 1333.3334
$ws.Range("K100").Value = 1182
$ws.Range("L100").Value = 1333.3334
$ws.Range("M100").Value = -641
$ws.Range("N100").Value = -2415.3334
$ws.Range("H138").Value = 3958.1016
$ws.Range("I138").Value = 985.1053000000001
$ws.Range("K138").Value = 2955.3159
$ws.Range("M138").Value = 2184.6841

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20837952
$ws.Range("I32").Value = 25002328
$ws.Range("K32").Value = 25002328
$ws.Range("M32").Value = -25002041
$ws.Range("H45").Value = 2400.6155
$ws.Range("I45").Value = 1666.5714
$ws.Range("J45").Value = 3257
$ws.Range("K45").Value = 1666.5714
$ws.Range("L45").Value = 3257
$ws.Range("M45").Value = -1289.5714
$ws.Range("N45").Value = -4011
$ws.Range("H74").Value = 12767.667
$ws.Range("J74").Value = 2000
$ws.Range("L74").Value = 2000
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 12767.667
$ws.Range("J77").Value = 2000
$ws.Range("L77").Value = 10000
$ws.Range("N77").Value = -18736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 250.66667
$ws.Range("I94").Value = 236.33333
$ws.Range("J94").Value = 265
$ws.Range("K94").Value = 236.33333
$ws.Range("L94").Value = 265
$ws.Range("M94").Value = 214.66667
$ws.Range("N94").Value = -1167
$ws.Range("H134").Value = 2637.4167
$ws.Range("I134").Value = 2704.4546
$ws.Range("K134").Value = 8113.3638
$ws.Range("M134").Value = -5578.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 6528.5
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 5371.3335
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 5371.3335
$ws.Range("M37").Value = -9893
$ws.Range("N37").Value = -5585.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7546.567
$ws.Range("J122").Value = 814.2857
$ws.Range("L122").Value = 7328.571300000001
$ws.Range("N122").Value = -12228.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.60714
$ws.Range("I2").Value = 52.81818
$ws.Range("J2").Value = 133.17647
$ws.Range("K2").Value = 52.81818
$ws.Range("L2").Value = 133.17647
$ws.Range("M2").Value = 60.18182
$ws.Range("N2").Value = -359.17647
$ws.Range("H33").Value = 1252250
$ws.Range("J33").Value = 1252250
$ws.Range("L33").Value = 1252250
$ws.Range("N33").Value = -1252754
$ws.Range("H44").Value = 4000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 4000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 4000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -5192
$ws.Range("H48").Value = 15000
$ws.Range("H80").Value = 2930.4119
$ws.Range("I80").Value = 2750.4167
$ws.Range("J80").Value = 3362.4
$ws.Range("K80").Value = 2750.4167
$ws.Range("L80").Value = 3362.4
$ws.Range("M80").Value = -1752.4167
$ws.Range("N80").Value = -5358.4
$ws.Range("H83").Value = 2930.4119
$ws.Range("I83").Value = 2750.4167
$ws.Range("J83").Value = 3362.4
$ws.Range("K83").Value = 13752.0835
$ws.Range("L83").Value = 16812
$ws.Range("M83").Value = -8760.083500000001
$ws.Range("N83").Value = -26796

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1300.75
$ws.Range("I16").Value = 1450.5
$ws.Range("J16").Value = 1151
$ws.Range("K16").Value = 1450.5
$ws.Range("L16").Value = 1151
$ws.Range("M16").Value = -1280.5
$ws.Range("N16").Value = -1491
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -22246
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -71232
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7808.5
$ws.Range("J41").Value = 7808.5
$ws.Range("L41").Value = 7808.5
$ws.Range("N41").Value = -8588.5
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 490.14285
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 445.25
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 1335.75
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -5675.75

Write-Host "Applied all Shinryu_Profits updates"